# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp footer
# - A handful of provinces tied in "Casos totales" and swapped rank order,
#   so their names need to move to the row that now matches their numbers
# - Refreshed case counts (Casos totales / Casos activos / Recuperados / Muertes)
#   for many provinces

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp update
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 14:22"

# Province-name realignments (ties in Casos totales reshuffled the ranking)
$ws.Range("A21").Value = "Sevilla"
$ws.Range("A22").Value = "Asturias"
$ws.Range("A25").Value = "Granada"
$ws.Range("A26").Value = "Pontevedra"
$ws.Range("A27").Value = "Cantabria"
$ws.Range("A38").Value = "Cadiz"
$ws.Range("A39").Value = "Castello/Castellon"
$ws.Range("A47").Value = "Gran Canaria"
$ws.Range("A48").Value = "Huesca"

# Refreshed numeric data
$ws.Range("B19").Value = 1932
$ws.Range("C19").Value = 289
$ws.Range("D19").Value = 1518
$ws.Range("E19").Value = 125

$ws.Range("B21").Value = 1713
$ws.Range("C21").Value = 105
$ws.Range("D21").Value = 1496
$ws.Range("E21").Value = 112

$ws.Range("B22").Value = 1679
$ws.Range("C22").Value = 244
$ws.Range("D22").Value = 1339
$ws.Range("E22").Value = 96

$ws.Range("B25").Value = 1550
$ws.Range("C25").Value = 182
$ws.Range("D25").Value = 1240
$ws.Range("E25").Value = 128

$ws.Range("B26").Value = 1536
$ws.Range("C26").Value = 333
$ws.Range("D26").Value = 1411
$ws.Range("E26").Value = 30

$ws.Range("D27").Value = 1288
$ws.Range("E27").Value = 84

$ws.Range("B32").Value = 1122
$ws.Range("D32").Value = 1444
$ws.Range("E32").Value = 59

$ws.Range("B33").Value = 1055
$ws.Range("D33").Value = 932
$ws.Range("E33").Value = 39

$ws.Range("B35").Value = 973
$ws.Range("C35").Value = 62
$ws.Range("D35").Value = 858
$ws.Range("E35").Value = 53

$ws.Range("B38").Value = 881
$ws.Range("C38").Value = 109
$ws.Range("D38").Value = 740
$ws.Range("E38").Value = 32

$ws.Range("B39").Value = 876
$ws.Range("C39").Value = 107
$ws.Range("D39").Value = 691
$ws.Range("E39").Value = 78

$ws.Range("B47").Value = 444
$ws.Range("C47").Value = 192
$ws.Range("D47").Value = 1444
$ws.Range("E47").Value = 25

$ws.Range("B48").Value = 419
$ws.Range("C48").Value = 62
$ws.Range("D48").Value = 315
$ws.Range("E48").Value = 42

$ws.Range("C50").Value = 50
$ws.Range("D50").Value = 298

$ws.Range("B52").Value = 292
$ws.Range("C52").Value = 19
$ws.Range("D52").Value = 257

$ws.Range("B56").Value = 67
$ws.Range("D56").Value = 1444

$ws.Range("B57").Value = 61
$ws.Range("D57").Value = 1444

$ws.Range("D59").Value = 1444
$ws.Range("D62").Value = 1444
$ws.Range("D64").Value = 1444
